$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update contribution percentages (Arbeitnehmer-/Arbeitgeberbeitrag RV): 9.8 -> 9.3
$ws.Range("B2").Value = 9.3
$ws.Range("B3").Value = 9.3

# Update Beitragsbemessungsgrenze RV Ost / West values
$ws.Range("B4").Value = 90600
$ws.Range("B5").Value = 89400

# Update Eintragsdatum (stored as text) 15.12.2023 -> 01.01.2024
$ws.Range("B6").Value = "01.01.2024"

# Move the active selection from A7 to B7
$ws.Range("B7").Select()
